# Generate Report for Handback
# The handback transform for the c25677a0-7d8b-4911-a3d4-bd24aee5419f file failed
# (handback file name did not match the handoff file name it should correspond to).
# Update the generated status report accordingly: the "Status" cells for that file
# flip from "Ready for handoff" to "Handback transform failed", and an explanatory
# message is recorded in the "Error Detail" column for each locale sheet.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$errorDetailZhCn = "Handback file name: wql4igiz.vle is different with handoff file name: c25677a0-7d8b-4911-a3d4-bd24aee5419f.4fdef314e5883d21b51c83fd7208556939a6b648.zh-cn."
$errorDetailDeDe = "Handback file name: wql4igiz.vle is different with handoff file name: c25677a0-7d8b-4911-a3d4-bd24aee5419f.4fdef314e5883d21b51c83fd7208556939a6b648.de-de."

# --- Overview sheet: status column for both locales, row for c25677a0 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# --- zh-cn detail sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("L3").Value = $errorDetailZhCn

# --- de-de detail sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("L3").Value = $errorDetailDeDe
